$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values get permuted across rows 2-13
$cols = @("A","B","D","E","F","G","H","Q","R")

# Snapshot the current ("before") values for each row so the shuffle can be
# applied without clobbering source data we still need to read.
$snapshot = @{}
for ($r = 2; $r -le 13; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Destination row -> source row mapping (after[row] = before[mapping[row]])
$mapping = @{
    2  = 6
    3  = 8
    4  = 11
    5  = 3
    6  = 10
    7  = 5
    8  = 2
    9  = 13
    10 = 7
    11 = 12
    12 = 9
    13 = 4
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
